# Update the "Estado de Cuenta" data table (rows 16-21) with the new
# periods (most recent first) and refreshed Valor Mora / Salario Basico
# figures, as described in the commit:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2307", "2306", "2305", "2304", "2301", "2212")
$valorMora = @(34666, 46400, 46400, 46400, 40000, 40000)
$salarioBasico = @(1160000, 1160000, 1160000, 1160000, 1160000, 1160000)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico[$i]
}
